$d = $word.ActiveDocument

# --- Table row height changes -------------------------------------------------
# Table 1, row 2 (the "very long title" row): atLeast 640 twips (32pt) -> 720 twips (36pt)
$t = $d.Tables.Item(1)
$row1 = $t.Rows.Item(2)
$row1.Height = 36

# Table 1, row 15 (the big spacer row): exact 11740 twips (587pt) -> 11660 twips (583pt)
$row2 = $t.Rows.Item(15)
$row2.Height = 583

# --- Font changes: Helvetica -> Arial for Sans_Bold / Sans_Normal styles -------
# (EMPTY_CELL_STYLE inherits its font from Sans_Normal and must stay untouched,
#  i.e. without an explicit w:rFonts of its own, so match on the style name
#  rather than the resolved/inherited Font.Name.)
$styles = $d.Styles
for ($i = 1; $i -le $styles.Count; $i++) {
    $s = $styles.Item($i)
    if ($s.NameLocal -eq "Sans_Bold" -or $s.NameLocal -eq "Sans_Normal") {
        $s.Font.NameAscii = "Arial"
        $s.Font.NameFarEast = "Arial"
        $s.Font.NameOther = "Arial"
        $s.Font.NameBi = "Arial"
    }
}
